$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 668.5294
$ws.Range("I9").Value = 358.53845
$ws.Range("K9").Value = 358.53845
$ws.Range("M9").Value = -189.53845
$ws.Range("H17").Value = 3736751
$ws.Range("J17").Value = 3736751
$ws.Range("L17").Value = 11210253
$ws.Range("N17").Value = -11210589
$ws.Range("H55").Value = 158.55556
$ws.Range("J55").Value = 224
$ws.Range("L55").Value = 224
$ws.Range("N55").Value = -652
$ws.Range("H58").Value = 1450
$ws.Range("I58").Value = 228.57143
$ws.Range("J58").Value = 10000
$ws.Range("K58").Value = 685.71429
$ws.Range("L58").Value = 30000
$ws.Range("M58").Value = -535.71429
$ws.Range("N58").Value = -30300
$ws.Range("H80").Value = 3034.7334
$ws.Range("I80").Value = 1609.5
$ws.Range("J80").Value = 3984.889
$ws.Range("K80").Value = 4828.5
$ws.Range("L80").Value = 11954.667
$ws.Range("M80").Value = -3830.5
$ws.Range("N80").Value = -13950.667
$ws.Range("H83").Value = 3034.7334
$ws.Range("I83").Value = 1609.5
$ws.Range("J83").Value = 3984.889
$ws.Range("K83").Value = 14485.5
$ws.Range("L83").Value = 35864.001
$ws.Range("M83").Value = -9493.5
$ws.Range("N83").Value = -45848.001
$ws.Range("H88").Value = 6545.8184
$ws.Range("I88").Value = 5400.2
$ws.Range("J88").Value = 7500.5
$ws.Range("K88").Value = 5400.2
$ws.Range("L88").Value = 7500.5
$ws.Range("M88").Value = -4994.2
$ws.Range("N88").Value = -8312.5
$ws.Range("H91").Value = 6545.8184
$ws.Range("I91").Value = 5400.2
$ws.Range("J91").Value = 7500.5
$ws.Range("K91").Value = 5400.2
$ws.Range("L91").Value = 7500.5
$ws.Range("M91").Value = -3996.2
$ws.Range("N91").Value = -10308.5
$ws.Range("H137").Value = 5345.2544
$ws.Range("I137").Value = 3109.9333
$ws.Range("J137").Value = 7657.6553
$ws.Range("K137").Value = 9329.7999
$ws.Range("L137").Value = 22972.9659
$ws.Range("M137").Value = -6779.7999
$ws.Range("N137").Value = -28072.9659

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 955.14
$ws.Range("I32").Value = 905.6383
$ws.Range("K32").Value = 905.6383
$ws.Range("M32").Value = -618.6383
$ws.Range("H54").Value = 39495
$ws.Range("J54").Value = 39495
$ws.Range("L54").Value = 39495
$ws.Range("N54").Value = -41033
$ws.Range("H61").Value = 20002824
$ws.Range("I61").Value = 21279132
$ws.Range("K61").Value = 21279132
$ws.Range("M61").Value = -21278920
$ws.Range("H122").Value = 1587.1
$ws.Range("I122").Value = 1046.4375
$ws.Range("K122").Value = 3139.3125
$ws.Range("M122").Value = -689.3125
$ws.Range("H132").Value = 27849536
$ws.Range("I132").Value = 6181.815
$ws.Range("K132").Value = 18545.445
$ws.Range("M132").Value = -16015.445
$ws.Range("H136").Value = 20002824
$ws.Range("I136").Value = 21279132
$ws.Range("K136").Value = 63837396
$ws.Range("M136").Value = -63834846

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2634138
$ws.Range("I134").Value = 2858918
$ws.Range("K134").Value = 8576754
$ws.Range("M134").Value = -8574219
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("N140").Value = 0
$ws.Range("L140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27029974
$ws.Range("I31").Value = 2029.9231
$ws.Range("K31").Value = 2029.9231
$ws.Range("M31").Value = -1734.9231
$ws.Range("H34").Value = 27029974
$ws.Range("I34").Value = 2029.9231
$ws.Range("K34").Value = 2029.9231
$ws.Range("M34").Value = -1827.9231
$ws.Range("H58").Value = 6735.7856
$ws.Range("I58").Value = 6804.636
$ws.Range("K58").Value = 6804.636
$ws.Range("M58").Value = -6601.636
$ws.Range("H62").Value = 4000.2856
$ws.Range("I62").Value = 3665.6667
$ws.Range("J62").Value = 4251.25
$ws.Range("K62").Value = 3665.6667
$ws.Range("L62").Value = 4251.25
$ws.Range("M62").Value = -3041.6667
$ws.Range("N62").Value = -5499.25
$ws.Range("H65").Value = 4000.2856
$ws.Range("I65").Value = 3665.6667
$ws.Range("J65").Value = 4251.25
$ws.Range("K65").Value = 18328.3335
$ws.Range("L65").Value = 21256.25
$ws.Range("M65").Value = -15208.3335
$ws.Range("N65").Value = -27496.25
$ws.Range("H86").Value = 9282.333000000001
$ws.Range("I86").Value = 9361.4
$ws.Range("J86").Value = 8887
$ws.Range("K86").Value = 9361.4
$ws.Range("L86").Value = 8887
$ws.Range("M86").Value = -8238.4
$ws.Range("N86").Value = -11133
$ws.Range("H89").Value = 9282.333000000001
$ws.Range("I89").Value = 9361.4
$ws.Range("J89").Value = 8887
$ws.Range("K89").Value = 46807
$ws.Range("L89").Value = 44435
$ws.Range("M89").Value = -41191
$ws.Range("N89").Value = -55667
$ws.Range("H105").Value = 28250
$ws.Range("I105").Value = 28250
$ws.Range("K105").Value = 28250
$ws.Range("M105").Value = -26503
$ws.Range("H132").Value = 48800.773
$ws.Range("I132").Value = 62876.79
$ws.Range("K132").Value = 188630.37
$ws.Range("M132").Value = -186100.37
$ws.Range("H134").Value = 2371.1316
$ws.Range("I134").Value = 2246.8386
$ws.Range("K134").Value = 6740.5158
$ws.Range("M134").Value = -4205.5158
$ws.Range("H136").Value = 6735.7856
$ws.Range("I136").Value = 6804.636
$ws.Range("K136").Value = 20413.908
$ws.Range("M136").Value = -17863.908
$ws.Range("H141").Value = 115658.89
$ws.Range("J141").Value = 139418.58
$ws.Range("L141").Value = 139418.58
$ws.Range("N141").Value = -149778.58

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J70").Value = 100
$ws.Range("L70").Value = 300
$ws.Range("N70").Value = -930
$ws.Range("J73").Value = 100
$ws.Range("L73").Value = 300
$ws.Range("N73").Value = -2484
$ws.Range("H128").Value = 99447.5
$ws.Range("I128").Value = 99447.5
$ws.Range("K128").Value = 298342.5
$ws.Range("M128").Value = -293362.5
$ws.Range("H134").Value = 1702.3438
$ws.Range("I134").Value = 1112.0968
$ws.Range("J134").Value = 20000
$ws.Range("K134").Value = 3336.2904
$ws.Range("L134").Value = 60000
$ws.Range("M134").Value = 1733.7096
$ws.Range("N134").Value = -70140
$ws.Range("H138").Value = 2111.4546
$ws.Range("I138").Value = 1770
$ws.Range("J138").Value = 3022
$ws.Range("K138").Value = 5310
$ws.Range("L138").Value = 9066
$ws.Range("M138").Value = -170
$ws.Range("N138").Value = -19346

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2550
$ws.Range("I31").Value = 2550
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2550
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = -2258
$ws.Range("H37").Value = 2550
$ws.Range("I37").Value = 2550
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 2550
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = -2273
$ws.Range("H49").Value = 27332.834
$ws.Range("H122").Value = 2769.6667
$ws.Range("I122").Value = 3182.5
$ws.Range("J122").Value = 2439.4
$ws.Range("K122").Value = 9547.5
$ws.Range("L122").Value = 7318.200000000001
$ws.Range("M122").Value = -7097.5
$ws.Range("N122").Value = -12218.2
$ws.Range("H126").Value = 37627812
$ws.Range("I126").Value = 20202500
$ws.Range("K126").Value = 60607500
$ws.Range("M126").Value = -60605030
$ws.Range("M31").ClearContents()
$ws.Range("M37").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4989.8887
$ws.Range("I7").Value = 4500
$ws.Range("J7").Value = 5381.8
$ws.Range("K7").Value = 4500
$ws.Range("L7").Value = 5381.8
$ws.Range("M7").Value = -4388
$ws.Range("N7").Value = -5605.8
$ws.Range("H22").Value = 2375.9412
$ws.Range("I22").Value = 1846.2222
$ws.Range("J22").Value = 2971.875
$ws.Range("K22").Value = 1846.2222
$ws.Range("L22").Value = 2971.875
$ws.Range("M22").Value = -1551.2222
$ws.Range("N22").Value = -3561.875
$ws.Range("H27").Value = 2375.9412
$ws.Range("I27").Value = 1846.2222
$ws.Range("J27").Value = 2971.875
$ws.Range("K27").Value = 1846.2222
$ws.Range("L27").Value = 2971.875
$ws.Range("M27").Value = -1739.2222
$ws.Range("N27").Value = -3185.875
$ws.Range("H40").Value = 3876.2144
$ws.Range("I40").Value = 3438.2856
$ws.Range("K40").Value = 3438.2856
$ws.Range("M40").Value = -3302.2856
$ws.Range("H100").Value = 2927.147
$ws.Range("I100").Value = 2219.7693
$ws.Range("K100").Value = 2219.7693
$ws.Range("M100").Value = -1678.7693
$ws.Range("H122").Value = 3238.5293
$ws.Range("I122").Value = 2567.5417
$ws.Range("J122").Value = 4848.9
$ws.Range("K122").Value = 7702.625100000001
$ws.Range("L122").Value = 14546.7
$ws.Range("M122").Value = -5252.625100000001
$ws.Range("N122").Value = -19446.7
$ws.Range("H126").Value = 4989.8887
$ws.Range("I126").Value = 4500
$ws.Range("J126").Value = 5381.8
$ws.Range("K126").Value = 13500
$ws.Range("L126").Value = 16145.4
$ws.Range("M126").Value = -11030
$ws.Range("N126").Value = -21085.4
$ws.Range("H132").Value = 34490696
$ws.Range("I132").Value = 4149.5347
$ws.Range("K132").Value = 12448.6041
$ws.Range("M132").Value = -9918.6041

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2170.3
$ws.Range("I81").Value = 1968.2142
$ws.Range("K81").Value = 3936.4284
$ws.Range("M81").Value = -2875.4284
$ws.Range("H84").Value = 2170.3
$ws.Range("I84").Value = 1968.2142
$ws.Range("K84").Value = 19682.142
$ws.Range("M84").Value = -14378.142
$ws.Range("H132").Value = 2179.15
$ws.Range("I132").Value = 2004.6389
$ws.Range("J132").Value = 3749.75
$ws.Range("K132").Value = 6013.9167
$ws.Range("L132").Value = 11249.25
$ws.Range("M132").Value = -3483.9167
$ws.Range("N132").Value = -16309.25
